# Add three new bullet paragraphs at the end of the document, after the
# paragraph ending in "...e mudando a fonte do título.", mirroring the
# existing "Aula 3" list structure (ilvl 1 = topic, ilvl 2 = sub-point).

$d = $word.ActiveDocument

# --- locate the very last paragraph in the document ---------------------
$lastIndex = $d.Paragraphs.Count
$lastPara  = $d.Paragraphs.Item($lastIndex)
$lastRng   = $lastPara.Range
$lastRng.Collapse(0)

# --- paragraph 1: "Criando classes para formatação de texto:" ----------
# (ilvl = 1, i.e. Word ListLevelNumber 2 -- same level as "Formatando o
#  título e subtítulo:")
$lastRng.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p1.Range.ListFormat.ListLevelNumber = 2
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertAfter(" Criando classes para formatação de texto:")

# --- paragraph 2: "Tomar cuidado ao usar as classes de cores..." -------
# (ilvl = 2, i.e. Word ListLevelNumber 3)
$r1end = $p1.Range
$r1end.Collapse(0)
$r1end.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2.Range.ListFormat.ListLevelNumber = 3
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertAfter("Tomar cuidado ao usar as classes de cores do bootstrap pois elas são semânticas!!!")

# --- paragraph 3: "Podemos criar classes customizadas..." --------------
# (ilvl = 2, i.e. Word ListLevelNumber 3)
$r2end = $p2.Range
$r2end.Collapse(0)
$r2end.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p3.Range.ListFormat.ListLevelNumber = 3
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter("Podemos criar classes customizadas para tudo o que o bootstrap não consegue fazer pra gente sem problema nenhum.")
